# Updated symbol list on Wed Jan  4 15:34:47 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns of the crypto symbol
# table on Sheet1 with the latest scraped quotes. Values are kept as
# literal text (leading apostrophe forces text entry so things like
# trailing zeros and "%" suffixes survive), and the style is reset to
# "Normal" right after so we do not leave a quote-prefix/text number
# format behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'255.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.19%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'28.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-3.78%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.204"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.91%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05857"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.52%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.690"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.16%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8709"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.74%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9580"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'11.53%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'2.57%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07121"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.49%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.03179"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.24%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'-1.03%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.001540"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.74%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.01056"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'5.02%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005967"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.31%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.496"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.42%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.212"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.79%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.225"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.3172"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.01%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.03490"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'4.58%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1308"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.40%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.528"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.23%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04190"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.65%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1366"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'2.73%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001224"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.35%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004553"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'9.93%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-0.12%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001466"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'1.06%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03816"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.09%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.005638"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.55%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1100"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'3.31%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002344"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.10%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009709"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.13%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005369"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.40%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.13%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.09497"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'5.55%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002130"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-12.66%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.13%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.13%"
$ws.Range("E50").Style = "Normal"
